# "statistics page is added" - Puanlama.xlsx edit
#
# The underlying data only kept the "McDonalds"/"BurgerKing" scores for a
# single person (Mustafa Tikir) instead of two people (Burak Simsek and
# Mustafa Tikir), and that person's scores were bumped to 7/7. The now-empty
# trailing rows are removed, which shrinks the table/used range from
# A1:D5 down to A1:D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 & 3 ("Burak Simsek") become "Mustafa Tikir", with Puan = 7 for both
$ws.Cells.Item(2, 1).Value = "Mustafa Tikir"
$ws.Cells.Item(3, 1).Value = "Mustafa Tikir"
$ws.Cells.Item(2, 4).Value = 7
$ws.Cells.Item(3, 4).Value = 7

# Rows 4 & 5 (the old "Mustafa Tikir" rows) are no longer needed - remove
# them bottom-up so row numbers of earlier rows don't shift underneath us.
# This also auto-shrinks the Table1 ListObject range/autofilter and the
# sheet dimension from A1:D5 to A1:D3.
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(4).Delete()

# Column A narrows slightly (13.140625 -> 13 characters of width)
$ws.Columns.Item(1).ColumnWidth = 12.14

# Selection moves from E3 to D3
$ws.Range("D3").Select()

# Window chrome size shrinks
$excel.ActiveWindow.Width = 17970
$excel.ActiveWindow.Height = 6075
